# android ws port - 30002
$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item(3)

# Populate the sheet with its new header/data (order chosen to match the
# original authoring sequence for the shared-strings table)
$wsInfo.Range("B1").Value = "_C_NAME"
$wsInfo.Range("C1").Value = "_C_HTML"
$wsInfo.Range("D1").Value = "_C_UDP_PORT"
$wsInfo.Range("E1").Value = "_C_MACADDR"

$wsInfo.Range("A2").Value = "_CUSTOM"
$wsInfo.Range("B2").Value = "CUSTOM"
$wsInfo.Range("C2").Value = "http://iotc365.com/test_locale/socket"
$wsInfo.Range("D2").Value = 11001
$wsInfo.Range("E2").Value = "00:08:22:d0:15:fc"

$wsInfo.Range("A3").Value = "_DISC"
$wsInfo.Range("B3").Value = "WEBDISC"
$wsInfo.Range("C3").Value = "http://static.iotc365.cn/socket"
$wsInfo.Range("D3").Value = 11002

# Rename the 3rd sheet ("Sheet3") to "web_url_info" and stamp the same text
# into A1
$wsInfo.Name = "web_url_info"
$wsInfo.Range("A1").Value = "web_url_info"

$wsInfo.Range("E3").Value = "00:08:22:d0:15:fc1"

# Set column widths to match the authored (best-fit) layout
$wsInfo.Columns.Item(1).ColumnWidth = 11.660714285714286
$wsInfo.Columns.Item(2).ColumnWidth = 8.910714285714286
$wsInfo.Columns.Item(3).ColumnWidth = 34.410714285714285
$wsInfo.Columns.Item(5).ColumnWidth = 14.785714285714286

# Make "web_url_info" the active/selected tab (clears tabSelected on "web1"
# and updates the workbook's activeTab index), then move the selection
$wsInfo.Activate()
$wsInfo.Range("E6").Select()
